$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2172.875
$ws.Range("J17").Value = 2172.875
$ws.Range("L17").Value = 6518.625
$ws.Range("N17").Value = -6854.625
$ws.Range("H33").Value = 739.1818
$ws.Range("I33").Value = 763.1
$ws.Range("J33").Value = 500
$ws.Range("K33").Value = 763.1
$ws.Range("L33").Value = 500
$ws.Range("M33").Value = -534.1
$ws.Range("N33").Value = -958
$ws.Range("H76").Value = 5283.0713
$ws.Range("I76").Value = 5987.25
$ws.Range("J76").Value = 4344.1665
$ws.Range("K76").Value = 5987.25
$ws.Range("L76").Value = 4344.1665
$ws.Range("M76").Value = -5672.25
$ws.Range("N76").Value = -4974.1665
$ws.Range("H79").Value = 5283.0713
$ws.Range("I79").Value = 5987.25
$ws.Range("J79").Value = 4344.1665
$ws.Range("K79").Value = 5987.25
$ws.Range("L79").Value = 4344.1665
$ws.Range("M79").Value = -4895.25
$ws.Range("N79").Value = -6528.1665
$ws.Range("H116").Value = 3202.6924
$ws.Range("J116").Value = 3486.2
$ws.Range("L116").Value = 3486.2
$ws.Range("N116").Value = -10370.2
$ws.Range("H118").Value = 1860.5
$ws.Range("I118").Value = 1790.8
$ws.Range("J118").Value = 2209
$ws.Range("K118").Value = 5372.4
$ws.Range("L118").Value = 6627
$ws.Range("M118").Value = -3715.4
$ws.Range("N118").Value = -9941
$ws.Range("H137").Value = 3214.5
$ws.Range("I137").Value = 2997.611
$ws.Range("J137").Value = 3702.5
$ws.Range("K137").Value = 8992.832999999999
$ws.Range("L137").Value = 11107.5
$ws.Range("M137").Value = -6442.832999999999
$ws.Range("N137").Value = -16207.5
$ws.Range("H138").Value = 2692.5823
$ws.Range("I138").Value = 1344
$ws.Range("J138").Value = 3316.926
$ws.Range("K138").Value = 4032
$ws.Range("L138").Value = 9950.778
$ws.Range("M138").Value = 1108
$ws.Range("N138").Value = -20230.778

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H11").Value = 2032.4
$ws.Range("I11").Value = 40.75
$ws.Range("J11").Value = 9999
$ws.Range("K11").Value = 40.75
$ws.Range("L11").Value = 9999
$ws.Range("M11").Value = 103.25
$ws.Range("N11").Value = -10287
$ws.Range("H97").Value = 929
$ws.Range("I97").Value = 929
$ws.Range("K97").Value = 929
$ws.Range("M97").Value = -433

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1760.6316
$ws.Range("I20").Value = 1858.5
$ws.Range("K20").Value = 1858.5
$ws.Range("M20").Value = -1611.5
$ws.Range("H134").Value = 2671.1428
$ws.Range("I134").Value = 1348.1613
$ws.Range("K134").Value = 4044.4839
$ws.Range("M134").Value = -1509.4839

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 8638084
$ws.Range("I99").Value = 3491848
$ws.Range("J99").Value = 10534066
$ws.Range("K99").Value = 3491848
$ws.Range("L99").Value = 10534066
$ws.Range("M99").Value = -3490350
$ws.Range("N99").Value = -10537062
$ws.Range("H126").Value = 8638084
$ws.Range("I126").Value = 3491848
$ws.Range("J126").Value = 10534066
$ws.Range("K126").Value = 10475544
$ws.Range("L126").Value = 31602198
$ws.Range("M126").Value = -10473074
$ws.Range("N126").Value = -31607138
$ws.Range("H132").Value = 3709
$ws.Range("J132").Value = 7605.4
$ws.Range("L132").Value = 22816.2
$ws.Range("N132").Value = -27876.2

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 3865.875
$ws.Range("J131").Value = 3956.1333
$ws.Range("L131").Value = 11868.3999
$ws.Range("N131").Value = -21948.3999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 77485.47
$ws.Range("I80").Value = 141923.62
$ws.Range("J80").Value = 3841.8572
$ws.Range("K80").Value = 141923.62
$ws.Range("L80").Value = 3841.8572
$ws.Range("M80").Value = -140925.62
$ws.Range("N80").Value = -5837.8572
$ws.Range("H83").Value = 77485.47
$ws.Range("I83").Value = 141923.62
$ws.Range("J83").Value = 3841.8572
$ws.Range("K83").Value = 709618.1
$ws.Range("L83").Value = 19209.286
$ws.Range("M83").Value = -704626.1
$ws.Range("N83").Value = -29193.286
$ws.Range("H126").Value = 4604.8823
$ws.Range("J126").Value = 5536.8184
$ws.Range("L126").Value = 16610.4552
$ws.Range("N126").Value = -21550.4552
$ws.Range("H132").Value = 3680.9473
$ws.Range("I132").Value = 3642.5881
$ws.Range("J132").Value = 4007
$ws.Range("K132").Value = 10927.7643
$ws.Range("L132").Value = 12021
$ws.Range("M132").Value = -8397.764299999999
$ws.Range("N132").Value = -17081

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3423.348
$ws.Range("I7").Value = 1825.5385
$ws.Range("K7").Value = 1825.5385
$ws.Range("M7").Value = -1713.5385
$ws.Range("H64").Value = 100150
$ws.Range("J64").Value = 100150
$ws.Range("L64").Value = 100150
$ws.Range("N64").Value = -100600
$ws.Range("H67").Value = 100150
$ws.Range("J67").Value = 100150
$ws.Range("L67").Value = 100150
$ws.Range("N67").Value = -101710
$ws.Range("H68").Value = 5320.032
$ws.Range("I68").Value = 3690.3333
$ws.Range("K68").Value = 3690.3333
$ws.Range("M68").Value = -2941.3333
$ws.Range("H71").Value = 5320.032
$ws.Range("I71").Value = 3690.3333
$ws.Range("K71").Value = 18451.6665
$ws.Range("M71").Value = -14707.6665
$ws.Range("H121").Value = 50000
$ws.Range("J121").Value = 50000
$ws.Range("L121").Value = 50000
$ws.Range("N121").Value = -53494
$ws.Range("H126").Value = 3423.348
$ws.Range("I126").Value = 1825.5385
$ws.Range("K126").Value = 5476.6155
$ws.Range("M126").Value = -3006.6155
$ws.Range("H132").Value = 4353.212
$ws.Range("I132").Value = 3898.3809
$ws.Range("J132").Value = 5149.1665
$ws.Range("K132").Value = 11695.1427
$ws.Range("L132").Value = 15447.4995
$ws.Range("M132").Value = -9165.1427
$ws.Range("N132").Value = -20507.4995
$ws.Range("H136").Value = 4335.244
$ws.Range("I136").Value = 2208.85
$ws.Range("K136").Value = 6626.549999999999
$ws.Range("M136").Value = -4076.549999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3808.8
$ws.Range("J62").Value = 4048
$ws.Range("L62").Value = 4048
$ws.Range("N62").Value = -5296
$ws.Range("H65").Value = 3808.8
$ws.Range("J65").Value = 4048
$ws.Range("L65").Value = 20240
$ws.Range("N65").Value = -26480
$ws.Range("H122").Value = 4492.5386
$ws.Range("I122").Value = 2334.5
$ws.Range("J122").Value = 14357.857
$ws.Range("K122").Value = 7003.5
$ws.Range("L122").Value = 43073.571
$ws.Range("M122").Value = -4553.5
$ws.Range("N122").Value = -47973.571
$ws.Range("H126").Value = 1433.3214
$ws.Range("I126").Value = 1501.4166
$ws.Range("K126").Value = 4504.2498
$ws.Range("M126").Value = -2034.2498
$ws.Range("H132").Value = 2254.5088
$ws.Range("I132").Value = 1986.5094
$ws.Range("J132").Value = 5805.5
$ws.Range("K132").Value = 5959.5282
$ws.Range("L132").Value = 17416.5
$ws.Range("M132").Value = -3429.5282
$ws.Range("N132").Value = -22476.5
